# Updated symbol list on Wed Dec 28 05:57:54 UTC 2022 with GitHub Actions
#
# This refreshes the scraped "Price" column (D) with newer quotes and fixes
# a handful of "Best/Worst in 24h" labels in the "Volume(1h)" column (E).
# Two rows (42/43) also got their coin identity (Coin/Link/Volume label)
# swapped along with brand-new price values.
#
# NOTE: column D stores numeric-looking values as literal TEXT (inline
# strings) in the workbook, not as numbers. Assigning a bare numeric string
# via COM (e.g. "243.29") would make Excel coerce it into a real number and
# lose formatting/precision (e.g. "0.001544" -> 0.001544 -> re-serialized
# differently, trailing zeros like "243.29" staying fine but others such as
# "0.00008820" would lose the trailing zero). To force Excel to keep these
# as text we prefix the value with a leading apostrophe, which is the
# standard way to enter "numbers as text" through the Excel UI/object model.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D (Price) text-value updates ---------------------------------
$priceUpdates = @{
    'D2'  = '243.29'
    'D3'  = '23.58'
    'D4'  = '5.295'
    'D5'  = '0.05767'
    'D7'  = '3.337'
    'D8'  = '0.8122'
    'D9'  = '0.8774'
    'D10' = '0.1390'
    'D11' = '0.07325'
    'D12' = '0.03094'
    'D13' = '0.03062'
    'D14' = '0.09319'
    'D15' = '3.870'
    'D16' = '0.001544'
    'D17' = '0.04706'
    'D18' = '0.0006062'
    'D19' = '0.006180'
    'D20' = '0.001295'
    'D21' = '0.00008820'
    'D22' = '3.582'
    'D23' = '2.143'
    'D25' = '0.1317'
    'D27' = '0.004603'
    'D28' = '0.0002354'
    'D40' = '0.03769'
    'D41' = '0.006375'
    'D42' = '0.1053'
    'D43' = '0.002636'
    'D44' = '0.007619'
    'D45' = '0.00005485'
    'D47' = '0.5912'
    'D48' = '0.001850'
    'D49' = '0.00002104'
    'D50' = '0.0002004'
}

foreach ($addr in $priceUpdates.Keys) {
    # Leading apostrophe forces text entry so the cell keeps its exact
    # string representation (no float coercion / trailing-zero loss).
    $ws.Range($addr).Value = "'" + $priceUpdates[$addr]
}

# --- Row 42 / 43: coin identity swap (Coin / Link / Volume label) --------
$ws.Range('B42').Value = 'BKEXToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range('E42').Value = '41BKEXTokenBKK'

$ws.Range('B43').Value = 'CEJI'
$ws.Range('C43').Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range('E43').Value = '42CEJICEJIWorstin24h'

# --- Column E (Volume(1h)) label-only updates -----------------------------
$ws.Range('E21').Value = '20NitroExNTXBestin24h'
$ws.Range('E41').Value = '40KickTokenKICK'
$ws.Range('E48').Value = '47BOLOBOLO'
